# Next use of TRIMRANGE
# Adds a "SecondDrill" worksheet with a real-estate listings table and a
# second TRIMRANGE/DROP example (ROWS + VLOOKUP), plus promotes the
# exceljet.net URL already present in P4 on FirstDrill into a real hyperlink.

$wb = $excel.ActiveWorkbook
$firstDrill = $wb.Worksheets.Item("FirstDrill")

# ---------------------------------------------------------------------------
# 1. Create the new worksheet right after FirstDrill
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add([System.Type]::Missing, $firstDrill)
$ws.Name = "SecondDrill"

# ---------------------------------------------------------------------------
# 2. Header row (row 4) - order matters for shared-string table layout
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = "Address"
$ws.Range("D4").Value = "Price"
$ws.Range("E4").Value = "Beds"
$ws.Range("F4").Value = "Baths"
$ws.Range("J4").Value = "Rows"

# ---------------------------------------------------------------------------
# 3. Address (C) & Size (G) columns, row by row
# ---------------------------------------------------------------------------
$addresses = @(
    "5335 Lake Rd", "4349 Cedar Ln", "2573 Maple Ave", "3773 Lake Dr",
    "1659 Main Ave", "2680 Lake Ave", "2386 oak Rd", "2636 oak st",
    "2095 Hill Dr", "2441 Main Ave", "1632 Lake Dr", "3390 Park st",
    "4752 Maple Ln"
)
$sizes = @(
    "2,003 For sale", "1,207 For sale", "3,454 sold", "3,455 sold",
    "1,041 For sale", "1,535 For sale", "1,831 For sale", "2,646 sold",
    "2,871 Sale Pending", "2,389 For sale", "21.58 sale pending",
    "1229 For sale", "2605 For sale"
)
for ($i = 0; $i -lt 13; $i++) {
    $r = 5 + $i
    $ws.Cells.Item($r, 3).Value = $addresses[$i]
    $ws.Cells.Item($r, 7).Value = $sizes[$i]
}

# ---------------------------------------------------------------------------
# 4. Remaining headers (typed last, after the data columns above)
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = "Size (sf)"
$ws.Range("H4").Value = "Status"

# ---------------------------------------------------------------------------
# 5. Status (H) column, row by row
# ---------------------------------------------------------------------------
$statuses = @(
    "For sale", "For sale", "sold", "sold", "For sale", "For sale",
    "For sale", "sold", "Sale Pending", "For sale", "sale pending",
    "For sale", "For sale"
)
for ($i = 0; $i -lt 13; $i++) {
    $r = 5 + $i
    $ws.Cells.Item($r, 8).Value = $statuses[$i]
}

# ---------------------------------------------------------------------------
# 6. ID header (reuses existing shared string) + numeric columns
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "ID"

$prices = @(535680, 423910, 355068, 473365, 349562, 3887740, 628305, 790649, 783209, 229702, 341303, 250047, 6749740)
$beds   = @(3, 1, 4, 4, 1, 2, 2, 4, 4, 3, 3, 1, 4)
$baths  = @(1, 4, 4, 2, 2, 3, 1, 3, 2, 1, 1, 2)

for ($i = 0; $i -lt 13; $i++) {
    $r = 5 + $i
    $ws.Cells.Item($r, 2).Value = $i + 1
    $ws.Cells.Item($r, 4).Value = $prices[$i]
    $ws.Cells.Item($r, 5).Value = $beds[$i]
    if ($i -lt 12) {
        $ws.Cells.Item($r, 6).Value = $baths[$i]
    }
}

# ---------------------------------------------------------------------------
# 7. Small lookup demo in columns J:K
# ---------------------------------------------------------------------------
$ws.Range("J5").Value = "ID"
$ws.Range("K5").Value = 12
$ws.Range("J6").Value = "Price"
$ws.Range("K4").Formula = "=ROWS(DROP(TRIMRANGE(B:H),1))"
$ws.Range("K6").Formula = "=VLOOKUP(K5,DROP(TRIMRANGE(B:H),1),3,0)"

Write-Output "done"
